# feat: add 2022-Q3 data
#
# Inserts a brand-new "2022-Q3" worksheet right after "总计", pushing the
# existing quarter sheets (2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q2) down
# by one tab position (their own content is untouched), and updates the
# "总计" (summary) sheet with the new quarter's row plus the now-visible
# 2021-Q2 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (summary) roll-up sheet: insert the 2022-Q3 entry
#    at the top of the data and shift everything else down by one row,
#    which surfaces a brand new trailing row for 2021-Q2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 2.57

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 2
$summary.Range("D3").Value = 2.77

$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 4.39

$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 2
$summary.Range("D5").Value = 1.48

$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 2
$summary.Range("D6").Value = 5.95

$summary.Range("A7").Value = 5
$summary.Range("A7").Font.Bold = $true
$summary.Range("A7").HorizontalAlignment = -4108
$summary.Range("A7").VerticalAlignment = -4160
$summary.Range("A7").Borders.LineStyle = 1
$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 2
$summary.Range("D7").Value = 4.84

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计" (this is what
#    naturally pushes 2022-Q2 / 2022-Q1 / 2021-Q4 / 2021-Q3 / 2021-Q2 one
#    slot to the right -- their data is left completely untouched).
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $afterSheet)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160
$q3.Range("B1:H1").Borders.LineStyle = 1

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'008763"
$q3.Range("C2").Value = "天弘越南市场股票（QDII）A"
$q3.Range("D2").Value = "'20.44"
$q3.Range("E2").Value = "'90.19"
$q3.Range("F2").Value = "'7.26"
$q3.Range("G2").Value = "'1.4839"
$q3.Range("H2").Value = 1

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'008764"
$q3.Range("C3").Value = "天弘越南市场股票（QDII）C"
$q3.Range("D3").Value = "'15.02"
$q3.Range("E3").Value = "'90.19"
$q3.Range("F3").Value = "'7.26"
$q3.Range("G3").Value = "'1.0905"
$q3.Range("H3").Value = 1

$q3.Range("A2:A3").Font.Bold = $true
$q3.Range("A2:A3").HorizontalAlignment = -4108
$q3.Range("A2:A3").VerticalAlignment = -4160
$q3.Range("A2:A3").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 3) Keep "2021-Q2" (now the last/trailing tab) as the selected sheet,
#    matching the original workbook's active-tab state.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()

